$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "Als..." / "moechte ich" / "damit/weil/denn" text for the
#    existing user-story rows (4-17) with the revised wording, and append
#    two new user stories as rows 18-19.
# ---------------------------------------------------------------------------

$ws.Range("D4").Value = "eine Liste der registrieten User abrufen"
$ws.Range("E4").Value = "ich einen Überblick über alle Userkonten erhalte und diese verwalten kann."

$ws.Range("D5").Value = "eine spzifische UserID abrufen"
$ws.Range("E5").Value = "ich die Details des Users angezeigt bekomme."

$ws.Range("D6").Value = "einen User erstellen"
$ws.Range("E6").Value = "dieser User bestimme Funktionen oder Dienste der Anwendung erhält."

$ws.Range("D7").Value = "einen User mit UserID angeben"
$ws.Range("E7").Value = "ich einen User unter einer spzifischen URL bearbeiten oder, falls dieser nicht existiert, erstellen kann."

$ws.Range("D8").Value = "spezifische Userdaten aktualisieren"
$ws.Range("E8").Value = "ich Daten ändern kann, ohne diese alle erneute eingeben zu müssen."

$ws.Range("D9").Value = "ein spezifisches Userprofil löschen"
$ws.Range("E9").Value = "ich meine Daten vollständig aus der Anwendung entfernen kann."

$ws.Range("D10").Value = "eine Liste aller vorhandener Events abfragen"
$ws.Range("E10").Value = "ich einen Überblick der vorhandenen Events habe."

$ws.Range("D11").Value = "Informationen über ein spezifisches Event abrufen"
$ws.Range("E11").Value = "ich alle Details erfahren kann, um zu entscheiden, ob ich teilnehmen möchte."

$ws.Range("D12").Value = "eine Liste aller Teilnehmer eines bestimmten Events abrufen"
$ws.Range("E12").Value = "ich weiß wer am Event teilnimmt, und ich mich mit diesen Personen vernetzen und vorbereiten kann."

$ws.Range("D13").Value = "ein Event erstellen und zu veröffentlichen"
$ws.Range("E13").Value = "ich andere Nutzer einladen und sie über die Details informieren kann."

$ws.Range("D14").Value = "ein Event aktualisieren können"
$ws.Range("E14").Value = "ich Änderungen an den Details des Events vornehmen kann."

$ws.Range("D15").Value = "spezifische Details eines Events teilweise aktualiseren "
$ws.Range("E15").Value = "ich schnell auf Änderungen reagieren kann, ohne alle Event-Informationen erneut einreichen zu müssen."

$ws.Range("D16").Value = "ein Event löschen kann"
$ws.Range("E16").Value = "ich es aus dem Veranstaltungskalender entfernen kann, falls es abgesagt wird oder nicht stattfinden kann."

$ws.Range("D17").Value = "einen User direkt zu meinem Event hinzufügen"
$ws.Range("E17").Value = "ich sicherstellen kann, dass bestimme Personen als Teilnehmer meines Events registriert sind"

# New row 18
$ws.Range("B18").Value = 15
$ws.Range("C18").Value = "Nutzer"
$ws.Range("D18").Value = "einen User einen Teilnehmer aus der Teilnehmerliste eines Events entfernen"
$ws.Range("E18").Value = "ich die Teilnehmerlsite aktuell halten und auf Änderungen reagieren kann."

# New row 19
$ws.Range("B19").Value = 16
$ws.Range("C19").Value = "Nutzer"
$ws.Range("D19").Value = "ein Event nach dessen Beendigung bewerten können"
$ws.Range("E19").Value = "ich Feedbackzu diesem Event bekommen kann und zukünftigen Teilnehmern eine Orientierung bieten kann."

# ---------------------------------------------------------------------------
# 2. The banding (top-border / bottom-border) shifted down by one row once
#    the new wording pushed things around, so re-apply the existing look
#    (copy format only, leaving the freshly written values untouched) from
#    rows that already carry the desired look.
# ---------------------------------------------------------------------------

# Row 8 (D8:E8) loses its banding -> becomes plain (like B5, unstyled)
$ws.Range("B5").Copy() | Out-Null
$ws.Range("D8:E8").PasteSpecial(-4122) | Out-Null

# Row 9 (D9:E9) becomes bottom-border banding (like B8)
$ws.Range("B8").Copy() | Out-Null
$ws.Range("D9:E9").PasteSpecial(-4122) | Out-Null

# Row 10 (D10:E10) becomes top-border banding (like B9)
$ws.Range("B9").Copy() | Out-Null
$ws.Range("D10:E10").PasteSpecial(-4122) | Out-Null

# Row 14 (D14:E14) loses its banding -> becomes plain (like B11, unstyled)
$ws.Range("B11").Copy() | Out-Null
$ws.Range("D14:E14").PasteSpecial(-4122) | Out-Null

# Row 15 (D15:E15) becomes bottom-border banding (like B8)
$ws.Range("B8").Copy() | Out-Null
$ws.Range("D15:E15").PasteSpecial(-4122) | Out-Null

# Row 16 (D16:E16) becomes top-border banding (like B9)
$ws.Range("B9").Copy() | Out-Null
$ws.Range("D16:E16").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Grow the table / AutoFilter range so it covers the two new rows.
# ---------------------------------------------------------------------------

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B3:E19")) | Out-Null

# ---------------------------------------------------------------------------
# 4. Widen the "moechte ich" / "damit/weil/denn" columns to fit the longer
#    text, and move the active selection the way the author left it.
# ---------------------------------------------------------------------------

$ws.Columns.Item(4).ColumnWidth = 70.5
$ws.Columns.Item(5).ColumnWidth = 96.5

$ws.Range("E20").Select() | Out-Null
